$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Issue Tracking")

# Row 8 (issue #6): append a new update line to the GetWireless comments (D8) cell
$d8 = $ws.Range("D8").Value()
$ws.Range("D8").Value = $d8 + "`n" + "•11/01 Does not seem to be related to the firmware switch. but to a specific SIM. We are asking Gemalto to check this on the modem."

# Row 7 (issue #5): append a new update line to the Resolution/Plan (E7) cell
$e7 = $ws.Range("E7").Value()
$ws.Range("E7").Value = $e7 + "`n" + "•01/11 It will be great if we can get a result by the end of the week so that we can add it on the release candidate."

# Row 8 (issue #6): append a new update line to the Resolution/Plan (E8) cell
$e8 = $ws.Range("E8").Value()
$ws.Range("E8").Value = $e8 + "`n" + "•01/11 The issue does not seem to be related to the firmware switch, but to something specific on the SIM. We are asking Gemalto to check this on the module."

# Row 8 is now taller because of the added lines
$ws.Rows.Item(8).RowHeight = 255

# Row 10 (issue #8): status moved from "New firmware" to "Closed"
$ws.Range("G10").Value = "Closed"

# Update the active view: scroll/frozen pane top-left cell and selection
$ws.Activate()
$ws.Range("E9").Select()
